$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5495
$ws1.Range("F10").Value = 2448
$ws1.Range("F12").Value = 73
$ws1.Range("F13").Value = 67
$ws1.Range("F14").Value = 2296
$ws1.Range("F15").Value = 190

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5495
$ws4.Range("F12").Value = 2448
$ws4.Range("F14").Value = 73
$ws4.Range("F16").Value = 67
$ws4.Range("F17").Value = 2296
$ws4.Range("F18").Value = 190
